$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cell (A1) rich-text formatting ---
# Original: "Supplementary Table S4.15: " as a single run, continuing with
# the rest of the description text (left untouched).
# New: "Supplementary Table" stays as the base (bold, sz14) run, " " becomes
# bold green (RGB 00B050), and "S4.15: " becomes bold (default/automatic
# colour, which resolves the same as the sheet's theme text colour).
$title = $ws.Range("A1")
$title.Characters(20, 1).Font.Color = 5287936
$title.Characters(21, 7).Font.Color = 0

# --- Shared-string relabelling for the data rows ---
foreach ($addr in @("A3", "A4", "A5")) { $ws.Range($addr).Value = "MCL14-BTB" }
foreach ($addr in @("A6", "A7", "A8")) { $ws.Range($addr).Value = "MCL21-BTB" }
foreach ($addr in @("A9", "A10", "A11")) { $ws.Range($addr).Value = "OGR25-BTB" }

foreach ($addr in @("B3", "B6", "B9")) { $ws.Range($addr).Value = "13-gene set" }
foreach ($addr in @("B4", "B7", "B10")) { $ws.Range($addr).Value = "17-gene set" }
foreach ($addr in @("B5", "B8", "B11")) { $ws.Range($addr).Value = "30-gene set" }

# --- Row heights: 20.1 -> 24.95 for the data rows (2-11) ---
$ws.Range("A2:A11").RowHeight = 24.95

# --- Updated P / P.adj. statistics ---
$ws.Range("G6").Value = 0.0049
$ws.Range("H6").Value = 0.011025

$ws.Range("G7").Value = 0.0098
$ws.Range("H7").Value = 0.01764

$ws.Range("H8").Value = 0.00735

$ws.Range("G9").Value = 0.0154
$ws.Range("H9").Value = 0.0231

$ws.Range("G10").Value = 0.000149
$ws.Range("H10").Value = 0.0006705

$ws.Range("G11").NumberFormat = "0.00E+00"
$ws.Range("G11").Value = 0.000029
$ws.Range("H11").Value = 0.000261

# --- Selection state ---
$ws.Range("A1:K1").Select()
